$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J10").Value = 3
$ws.Range("M10").Value = "71.6 MPH"
$ws.Range("J11").Value = 1
$ws.Range("M12").Value = "22.4°"
$ws.Range("J14").Value = "Roblez"
$ws.Range("M14").Value = "Line Drive"
$ws.Range("M15").Value = "Single"
$ws.Range("J16").Value = "88-90 MPH"
$ws.Range("J17").Value = "CB,FB,CH"
$ws.Range("J19").Value = 7
$ws.Range("M19").Value = "20.44 MPH"
$ws.Range("M21").Value = "7.81°"
$ws.Range("J23").Value = "Plum"
$ws.Range("M23").Value = "Undefined"
$ws.Range("M24").Value = "Undefined"
$ws.Range("J25").Value = "84-86 MPH"
$ws.Range("J26").Value = "SL,FB,CH"
$ws.Range("J28").Value = 7
$ws.Range("M28").Value = "nan MPH"
$ws.Range("J29").Value = 0
$ws.Range("M30").Value = "nan°"
$ws.Range("J32").Value = "Plum"
$ws.Range("M32").Value = "Undefined"
$ws.Range("J33").Value = "Right"
$ws.Range("M33").Value = "Undefined"
$ws.Range("J34").Value = "84-86 MPH"
$ws.Range("J35").Value = "SL,FB,CH"
$ws.Range("J37").Value = 5
$ws.Range("M37").Value = "99.63 MPH"
$ws.Range("J38").Value = 2
$ws.Range("M39").Value = "10.11°"
$ws.Range("J41").Value = "Herbst"
$ws.Range("M41").Value = "Ground Ball"
$ws.Range("M42").Value = "Double"
$ws.Range("J43").Value = "83-85 MPH"
$ws.Range("J44").Value = "SL,CB,FB,CH"
$ws.Range("J46").Value = 9
$ws.Range("M46").Value = "53.02 MPH"
$ws.Range("M48").Value = "31.83°"
$ws.Range("J50").Value = "Thompson"
$ws.Range("M50").Value = "Popup"
$ws.Range("J51").Value = "Left"
$ws.Range("M51").Value = "Out"
$ws.Range("J52").Value = "84-84 MPH"
$ws.Range("J53").Value = "SL,FB,CH"
$ws.Range("J61").Value = 4
$ws.Range("M61").Value = "nan MPH"
$ws.Range("J62").Value = 0
$ws.Range("M63").Value = "nan°"
$ws.Range("J65").Value = "Roblez"
$ws.Range("J67").Value = "88-90 MPH"
$ws.Range("J68").Value = "CB,FB,CH"
